# Updated home text, baseline and followup
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1. Insert the two new rows after the current row 5
#    (new row 6: "reacting to online learning", new row 7: "teaching hours")
# ---------------------------------------------------------------------------
$ws.Range("A6:A7").EntireRow.Insert() | Out-Null

$ws.Range("A6").Value = "[If yes to children in K-12]`nIf  your child has already begun the 2020-2021 school year, how are they reacting to online learning?"
$ws.Range("B6").Value = "$bullet`tVery well`n$bullet`tSomewhat well`n$bullet`tWell`n$bullet`tNot well `n$bullet`tVery poorly"
$ws.Range("C6").Value = "Developed by RAPID Team"
$ws.Range("D6").Value = "Current 23"
$ws.Rows.Item(6).RowHeight = 85

$ws.Range("A7").Value = "[If yes to children in K-12]`nIncluding hours spent during weekdays and weekends, about how many hours did you spend on teaching activities with your school-aged child(ren) in this household during the last 7 days? Enter the total number of hours. If none, enter 0."
$ws.Range("B7").Value = "Open Response"
$ws.Range("C7").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = "Current 23"
$ws.Rows.Item(7).RowHeight = 68

# ---------------------------------------------------------------------------
# 2. Insert two more new rows after the (now shifted) "Will you use a child
#    care provider..." row, which is row 8 after the previous insert.
#    (new row 9: kindergarten wait, new row 10: if yes why)
# ---------------------------------------------------------------------------
$ws.Range("A9:A10").EntireRow.Insert() | Out-Null

$ws.Range("A9").Value = "If you have a child that was due to be entering kindergarten this fall, have you decided to wait until next fall (2021) instead because of the pandemic?"
$ws.Range("B9").Value = "$bullet`tYes`n$bullet`tNo `n$bullet`tNot applicable"
$ws.Range("C9").Value = "Developed by RAPID Team"
$ws.Range("D9").Value = "Current 23"
$ws.Rows.Item(9).RowHeight = 51

$ws.Range("A10").Value = "If yes, why? Select all that apply. "
$ws.Range("B10").Value = "$bullet`tSafety`n$bullet`tUncertain about the plan for school (in person/online)`n$bullet`tNot able to manage online instruction for my child along with my other responsibilities (work, etc.)"
$ws.Range("C10").Value = "Developed by RAPID Team"
$ws.Range("D10").Value = "Current 23"
$ws.Rows.Item(10).RowHeight = 119

# ---------------------------------------------------------------------------
# 3. Update the "Occurrence" column (D) text for all the pre-existing rows
#    from "Current\n21" to "Current\n21, 23" (rows 2-5 untouched by inserts,
#    plus the rows that have since shifted down: old 6->8, old7->11,
#    old8->12, old9->13, old10->14).
# ---------------------------------------------------------------------------
$currentText = "Current`n21, 23"
foreach ($r in 2,3,4,5,8,11,12,13,14) {
    $ws.Range("D$r").Value = $currentText
}

# ---------------------------------------------------------------------------
# 4. Append two brand-new rows at the end of the table (new rows 15 and 16).
#    Copy the style from the row above (row 14) first so the new rows pick
#    up the same wrap/vertical-top formatting used throughout the table.
# ---------------------------------------------------------------------------
$ws.Range("A14:D14").Copy() | Out-Null
$ws.Range("A15:D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend on face-to-face activities with your children 0-5 in this household during the last 7 days? Enter the total number of hours. If none, enter 0."
$ws.Range("B15").Value = "Open Response"
$ws.Range("C15").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("C15").WrapText = $true
$ws.Range("D15").Value = "Current 23"
$ws.Rows.Item(15).RowHeight = 51

$ws.Range("A16").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend working during the last 7 days? Enter the total number of hours. If none, enter 0 or NA if not working currently. "
$ws.Range("B16").Value = "Open Response"
$ws.Range("C16").Value = "RAPID Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("C16").WrapText = $true
$ws.Range("D16").Value = "Current 23"
$ws.Rows.Item(16).RowHeight = 51

# ---------------------------------------------------------------------------
# 6. Update the view to match: scrolled so row 13 is the top row, with B15
#    selected.
# ---------------------------------------------------------------------------
$ws.Range("B15").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
